# Give_Feedback.xlsx fix-up:
#  - Column A (rows 2-6) held the placeholder text "Test Account" in every
#    row. Replace it with the real numeric account id 443671 (a genuine
#    number, just displayed via a text-style left-aligned format).
#  - Leave the selection on A2 (first data cell in the id column) instead
#    of the old F4 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the numeric id first, while the cells are still in their default
# (General) format, so it is stored as a real number.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = 443671
}

# Now apply the display format (same text-number style already used by
# the membership-term columns L/M) without disturbing the stored value.
$idRange = $ws.Range("A2:A6")
$idRange.NumberFormat = "@"
$idRange.HorizontalAlignment = -4131  # xlLeft

$ws.Range("A2").Select() | Out-Null
